$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 12.35843671810672
$ws.Cells.Item(2, 3).Value = 5.508036341809289
$ws.Cells.Item(2, 4).Value = 5.999005911920817
$ws.Cells.Item(2, 5).Value = 16.35999277170716
$ws.Cells.Item(2, 7).Value = 41.71122937633841
$ws.Cells.Item(2, 8).Value = 17.06541419583024
$ws.Cells.Item(2, 11).Value = 11.63034214925159
$ws.Cells.Item(2, 14).Value = 19.86073858918607

$ws.Cells.Item(3, 2).Value = 12.04316319846023
$ws.Cells.Item(3, 3).Value = 5.162126300115098
$ws.Cells.Item(3, 4).Value = 5.883397141789339
$ws.Cells.Item(3, 5).Value = 15.43925682576197
$ws.Cells.Item(3, 7).Value = 41.32267618239997
$ws.Cells.Item(3, 8).Value = 17.0629749853773
$ws.Cells.Item(3, 11).Value = 11.40748985987621
$ws.Cells.Item(3, 14).Value = 19.91139020816327

$ws.Cells.Item(4, 2).Value = 11.84926231524708
$ws.Cells.Item(4, 3).Value = 4.936938499831816
$ws.Cells.Item(4, 4).Value = 5.813138927669579
$ws.Cells.Item(4, 5).Value = 14.85031651072761
$ws.Cells.Item(4, 7).Value = 41.09515039061506
$ws.Cells.Item(4, 8).Value = 17.06492943537852
$ws.Cells.Item(4, 11).Value = 11.27200816352311
$ws.Cells.Item(4, 14).Value = 19.94438800473879

$ws.Cells.Item(5, 2).Value = 11.77029817436059
$ws.Cells.Item(5, 3).Value = 4.841956118192806
$ws.Cells.Item(5, 4).Value = 5.784733618210627
$ws.Cells.Item(5, 5).Value = 14.60465718399173
$ws.Cells.Item(5, 7).Value = 41.00529488824296
$ws.Cells.Item(5, 8).Value = 17.0665927371599
$ws.Cells.Item(5, 11).Value = 11.21722003730775
$ws.Cells.Item(5, 14).Value = 19.95831195949065

$ws.Cells.Item(6, 2).Value = 11.75719335736249
$ws.Cells.Item(6, 3).Value = 4.82836110280795
$ws.Cells.Item(6, 4).Value = 5.780031813462654
$ws.Cells.Item(6, 5).Value = 14.56353276509722
$ws.Cells.Item(6, 7).Value = 40.99054973346959
$ws.Cells.Item(6, 8).Value = 17.06692123850242
$ws.Cells.Item(6, 11).Value = 11.2081503910277
$ws.Cells.Item(6, 14).Value = 19.96065283508931

$ws.Cells.Item(7, 2).Value = 11.84819698747286
$ws.Cells.Item(7, 3).Value = 4.935670553804009
$ws.Cells.Item(7, 4).Value = 5.812754873930466
$ws.Cells.Item(7, 5).Value = 14.84702598442441
$ws.Cells.Item(7, 7).Value = 41.09392686793665
$ws.Cells.Item(7, 8).Value = 17.06494835945247
$ws.Cells.Item(7, 11).Value = 11.27126745523602
$ws.Cells.Item(7, 14).Value = 19.94457385650033

$ws.Cells.Item(8, 2).Value = 12.24987539025468
$ws.Cells.Item(8, 3).Value = 5.391422816637833
$ws.Cells.Item(8, 4).Value = 5.959018990996423
$ws.Cells.Item(8, 5).Value = 16.04757105097803
$ws.Cells.Item(8, 7).Value = 41.57501637830266
$ws.Cells.Item(8, 8).Value = 17.06385605281451
$ws.Cells.Item(8, 11).Value = 11.55327152310485
$ws.Cells.Item(8, 14).Value = 19.87780940499683

$ws.Cells.Item(9, 2).Value = 13.02944471616903
$ws.Cells.Item(9, 3).Value = 6.183671058938591
$ws.Cells.Item(9, 4).Value = 6.249737708914597
$ws.Cells.Item(9, 5).Value = 18.24631831249159
$ws.Cells.Item(9, 7).Value = 42.60204768780628
$ws.Cells.Item(9, 8).Value = 17.08914389940609
$ws.Cells.Item(9, 11).Value = 12.11342626673942
$ws.Cells.Item(9, 14).Value = 19.761942631747

$ws.Cells.Item(10, 2).Value = 13.59024232828602
$ws.Cells.Item(10, 3).Value = 6.704101913143708
$ws.Cells.Item(10, 4).Value = 6.463264250959601
$ws.Cells.Item(10, 5).Value = 19.87580381553409
$ws.Cells.Item(10, 7).Value = 43.4018034186447
$ws.Cells.Item(10, 8).Value = 17.12447470300155
$ws.Cells.Item(10, 11).Value = 12.52466223366704
$ws.Cells.Item(10, 14).Value = 19.68599811417862

$ws.Cells.Item(11, 2).Value = 13.84143422811712
$ws.Cells.Item(11, 3).Value = 6.927553905406995
$ws.Cells.Item(11, 4).Value = 6.55989818150348
$ws.Cells.Item(11, 5).Value = 20.57566641787696
$ws.Cells.Item(11, 7).Value = 43.7741312013851
$ws.Cells.Item(11, 8).Value = 17.14418222197461
$ws.Cells.Item(11, 11).Value = 12.71073827541505
$ws.Cells.Item(11, 14).Value = 19.65344415347302

$ws.Cells.Item(12, 2).Value = 13.93589036150503
$ws.Cells.Item(12, 3).Value = 7.010264674658407
$ws.Cells.Item(12, 4).Value = 6.596380408425444
$ws.Cells.Item(12, 5).Value = 20.83477083854785
$ws.Cells.Item(12, 7).Value = 43.91623189579694
$ws.Cells.Item(12, 8).Value = 17.15216656711981
$ws.Cells.Item(12, 11).Value = 12.78098325885609
$ws.Cells.Item(12, 14).Value = 19.64140372351147

$ws.Cells.Item(13, 2).Value = 13.91557863760663
$ws.Cells.Item(13, 3).Value = 6.992536080183521
$ws.Cells.Item(13, 4).Value = 6.588528841853441
$ws.Cells.Item(13, 5).Value = 20.77923073570845
$ws.Cells.Item(13, 7).Value = 43.88558054596287
$ws.Cells.Item(13, 8).Value = 17.15042382468157
$ws.Cells.Item(13, 11).Value = 12.76586558561423
$ws.Cells.Item(13, 14).Value = 19.64398407265577

$ws.Cells.Item(14, 2).Value = 13.84921906622965
$ws.Cells.Item(14, 3).Value = 6.934396724746596
$ws.Cells.Item(14, 4).Value = 6.562902035972342
$ws.Cells.Item(14, 5).Value = 20.59710150912841
$ws.Cells.Item(14, 7).Value = 43.78580037139377
$ws.Cells.Item(14, 8).Value = 17.14482865523439
$ws.Cells.Item(14, 11).Value = 12.71652215446965
$ws.Cells.Item(14, 14).Value = 19.65244782492282

$ws.Cells.Item(15, 2).Value = 13.80848246736562
$ws.Cells.Item(15, 3).Value = 6.898536712283687
$ws.Cells.Item(15, 4).Value = 6.547189306028079
$ws.Cells.Item(15, 5).Value = 20.48477238064643
$ws.Cells.Item(15, 7).Value = 43.72482298311797
$ws.Cells.Item(15, 8).Value = 17.141469327146
$ws.Cells.Item(15, 11).Value = 12.68626734369738
$ws.Cells.Item(15, 14).Value = 19.65766951139263

$ws.Cells.Item(16, 2).Value = 13.57373865534603
$ws.Cells.Item(16, 3).Value = 6.689231782521393
$ws.Cells.Item(16, 4).Value = 6.456935549150668
$ws.Cells.Item(16, 5).Value = 19.82923607508833
$ws.Cells.Item(16, 7).Value = 43.37763237203367
$ws.Cells.Item(16, 8).Value = 17.12325982926185
$ws.Cells.Item(16, 11).Value = 12.5124751981324
$ws.Cells.Item(16, 14).Value = 19.68816573923691

$ws.Cells.Item(17, 2).Value = 13.4286527450242
$ws.Cells.Item(17, 3).Value = 6.557429013684259
$ws.Cells.Item(17, 4).Value = 6.401411269704591
$ws.Cells.Item(17, 5).Value = 19.4165073816253
$ws.Cells.Item(17, 7).Value = 43.16674097068108
$ws.Cells.Item(17, 8).Value = 17.11301927565504
$ws.Cells.Item(17, 11).Value = 12.40554967310865
$ws.Cells.Item(17, 14).Value = 19.70738504689623

$ws.Cells.Item(18, 2).Value = 13.34483995648607
$ws.Cells.Item(18, 3).Value = 6.480367776051851
$ws.Cells.Item(18, 4).Value = 6.369430288613438
$ws.Cells.Item(18, 5).Value = 19.17521806515202
$ws.Cells.Item(18, 7).Value = 43.04625037273693
$ws.Cells.Item(18, 8).Value = 17.10747146123082
$ws.Cells.Item(18, 11).Value = 12.34395930313338
$ws.Cells.Item(18, 14).Value = 19.71862710532761

$ws.Cells.Item(19, 2).Value = 13.31640322894468
$ws.Cells.Item(19, 3).Value = 6.454060998904955
$ws.Cells.Item(19, 4).Value = 6.358595589519626
$ws.Cells.Item(19, 5).Value = 19.0928509850112
$ws.Cells.Item(19, 7).Value = 43.00559670747041
$ws.Cells.Item(19, 8).Value = 17.10565188310145
$ws.Cells.Item(19, 11).Value = 12.32309280106537
$ws.Cells.Item(19, 14).Value = 19.72246568632487

$ws.Cells.Item(20, 2).Value = 13.44413573946574
$ws.Cells.Item(20, 3).Value = 6.571589154979307
$ws.Cells.Item(20, 4).Value = 6.407326858132123
$ws.Cells.Item(20, 5).Value = 19.46084630209838
$ws.Cells.Item(20, 7).Value = 43.18910787422461
$ws.Cells.Item(20, 8).Value = 17.11407398280989
$ws.Cells.Item(20, 11).Value = 12.41694190712438
$ws.Cells.Item(20, 14).Value = 19.70531969826701

$ws.Cells.Item(21, 2).Value = 13.86872929076358
$ws.Cells.Item(21, 3).Value = 6.951525329377013
$ws.Cells.Item(21, 4).Value = 6.570432561868424
$ws.Cells.Item(21, 5).Value = 20.65075762481569
$ws.Cells.Item(21, 7).Value = 43.81507908383972
$ws.Cells.Item(21, 8).Value = 17.14645794994036
$ws.Cells.Item(21, 11).Value = 12.73102200727683
$ws.Cells.Item(21, 14).Value = 19.64995402215929

$ws.Cells.Item(22, 2).Value = 14.14230955019987
$ws.Cells.Item(22, 3).Value = 7.188731301721742
$ws.Cells.Item(22, 4).Value = 6.676368256521434
$ws.Cells.Item(22, 5).Value = 21.39395270758272
$ws.Cells.Item(22, 7).Value = 44.23059230775594
$ws.Cells.Item(22, 8).Value = 17.17066174601726
$ws.Cells.Item(22, 11).Value = 12.93498916661796
$ws.Cells.Item(22, 14).Value = 19.61544268699685

$ws.Cells.Item(23, 2).Value = 13.99668434116548
$ws.Cells.Item(23, 3).Value = 7.06314385737388
$ws.Cells.Item(23, 4).Value = 6.619901431253409
$ws.Cells.Item(23, 5).Value = 21.00043908207999
$ws.Cells.Item(23, 7).Value = 44.00827779276539
$ws.Cells.Item(23, 8).Value = 17.15746618401437
$ws.Cells.Item(23, 11).Value = 12.82627084776438
$ws.Cells.Item(23, 14).Value = 19.63370879377146

$ws.Cells.Item(24, 2).Value = 13.43713711878438
$ws.Cells.Item(24, 3).Value = 6.565191361451826
$ws.Cells.Item(24, 4).Value = 6.404652602451957
$ws.Cells.Item(24, 5).Value = 19.4408131573359
$ws.Cells.Item(24, 7).Value = 43.17899344059558
$ws.Cells.Item(24, 8).Value = 17.11359609168967
$ws.Cells.Item(24, 11).Value = 12.41179183825663
$ws.Cells.Item(24, 14).Value = 19.70625284235848

$ws.Cells.Item(25, 2).Value = 12.82017054597568
$ws.Cells.Item(25, 3).Value = 5.980192960530091
$ws.Cells.Item(25, 4).Value = 6.170929981575554
$ws.Cells.Item(25, 5).Value = 17.64355665442167
$ws.Cells.Item(25, 7).Value = 42.31583547664321
$ws.Cells.Item(25, 8).Value = 17.07936140139032
$ws.Cells.Item(25, 11).Value = 11.96161512267082
$ws.Cells.Item(25, 14).Value = 19.79167539701578
